$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2026-02-08T18:08:24"
$ws.Range("V4").Value = 368
$ws.Range("W4").Value = 323.72
$ws.Range("Y4").Value = 210.31
$ws.Range("Z4").Value = 304.07
$ws.Range("V6").Value = -25.01
$ws.Range("W6").Value = -20.72
$ws.Range("Y6").Value = -10.73
$ws.Range("Z6").Value = -13.38
$ws.Range("V8").Value = 20.7
$ws.Range("V9").Value = 348.15
$ws.Range("W9").Value = 303.47
$ws.Range("X9").Value = 208.28
$ws.Range("Y9").Value = 206.19
$ws.Range("Z9").Value = 294.48
$ws.Range("V11").Value = -44.86
$ws.Range("W11").Value = -40.97
$ws.Range("X11").Value = -16.04
$ws.Range("Y11").Value = -14.85
$ws.Range("Z11").Value = -22.97
$ws.Range("V13").Value = 20.7
$ws.Range("V14").Value = 82
$ws.Range("W14").Value = 57.37
$ws.Range("X14").Value = 208.28
$ws.Range("Y14").Value = 206.19
$ws.Range("Z14").Value = 294.48
$ws.Range("V15").Value = -266.14
$ws.Range("W15").Value = -246.1
$ws.Range("V16").Value = -44.86
$ws.Range("W16").Value = -40.97
$ws.Range("X16").Value = -16.04
$ws.Range("Y16").Value = -14.85
$ws.Range("Z16").Value = -22.97
$ws.Range("V18").Value = 20.7
$ws.Range("V19").Value = 113.04
$ws.Range("W19").Value = 85.48
$ws.Range("X19").Value = 114.49
$ws.Range("Y19").Value = 114.65
$ws.Range("Z19").Value = 78.34
$ws.Range("V20").Value = -249.86
$ws.Range("W20").Value = -233.15
$ws.Range("Y20").Value = -92.90000000000001
$ws.Range("Z20").Value = -221.42
$ws.Range("V21").Value = -30.11
$ws.Range("W21").Value = -25.81
$ws.Range("X21").Value = -14.87
$ws.Range("Y21").Value = -13.49
$ws.Range("Z21").Value = -17.69
$ws.Range("V23").Value = 20.7
$ws.Range("V24").Value = 362.89
$ws.Range("W24").Value = 318.63
$ws.Range("X24").Value = 209.44
$ws.Range("Y24").Value = 207.55
$ws.Range("Z24").Value = 299.76
$ws.Range("V26").Value = -30.11
$ws.Range("W26").Value = -25.81
$ws.Range("X26").Value = -14.87
$ws.Range("Y26").Value = -13.49
$ws.Range("Z26").Value = -17.69
$ws.Range("V28").Value = 20.7
$ws.Range("V29").Value = 108.08
$ws.Range("W29").Value = 80.27
$ws.Range("Y29").Value = 111.39
$ws.Range("Z29").Value = 73.06
$ws.Range("V30").Value = -249.86
$ws.Range("W30").Value = -233.15
$ws.Range("Y30").Value = -92.90000000000001
$ws.Range("Z30").Value = -221.42
$ws.Range("V31").Value = -35.07
$ws.Range("W31").Value = -31.03
$ws.Range("Y31").Value = -16.75
$ws.Range("Z31").Value = -22.97
$ws.Range("V33").Value = 20.7
$ws.Range("V34").Value = 77.18000000000001
$ws.Range("V35").Value = -266.14
$ws.Range("W35").Value = -246.1
$ws.Range("V36").Value = -49.68
$ws.Range("W36").Value = -46.22
$ws.Range("X36").Value = -16.23
$ws.Range("Y36").Value = -15.04
$ws.Range("Z36").Value = -24.06
$ws.Range("W37").Value = -51.11
$ws.Range("X37").Value = -207.08
$ws.Range("Y37").Value = -205
$ws.Range("Z37").Value = -292.39
$ws.Range("V38").Value = 20.7
$ws.Range("V39").Value = 368
$ws.Range("W39").Value = 323.72
$ws.Range("Y39").Value = 210.31
$ws.Range("Z39").Value = 304.07
$ws.Range("V41").Value = -25.01
$ws.Range("W41").Value = -20.72
$ws.Range("Y41").Value = -10.73
$ws.Range("Z41").Value = -13.38
$ws.Range("V43").Value = 20.7
$ws.Range("V44").Value = 388.96
$ws.Range("W44").Value = 341.37
$ws.Range("Y44").Value = 219.28
$ws.Range("Z44").Value = 316.82
$ws.Range("V46").Value = -4.05
$ws.Range("W46").Value = -3.07
$ws.Range("Y46").Value = -1.75
$ws.Range("Z46").Value = -0.63
$ws.Range("V48").Value = 20.7
$ws.Range("V49").Value = 369.96
$ws.Range("W49").Value = 338.02
$ws.Range("Y49").Value = 229.06
$ws.Range("Z49").Value = 331.37
$ws.Range("V51").Value = -23.05
$ws.Range("W51").Value = -6.42
$ws.Range("Y51").Value = 8.02
$ws.Range("Z51").Value = 13.92
$ws.Range("V53").Value = 20.7
$ws.Range("V54").Value = 383.57
$ws.Range("W54").Value = 334.41
$ws.Range("X54").Value = 221.44
$ws.Range("Y54").Value = 221.7
$ws.Range("Z54").Value = 315.24
$ws.Range("V56").Value = -9.43
$ws.Range("W56").Value = -10.03
$ws.Range("X56").Value = -2.88
$ws.Range("Y56").Value = 0.67
$ws.Range("Z56").Value = -2.21
$ws.Range("V58").Value = 20.7
$ws.Range("V59").Value = 400.99
$ws.Range("W59").Value = 351.83
$ws.Range("Y59").Value = 224.86
$ws.Range("Z59").Value = 324.59
$ws.Range("V61").Value = 7.99
$ws.Range("W61").Value = 7.39
$ws.Range("Y61").Value = 3.82
$ws.Range("Z61").Value = 7.14
$ws.Range("V63").Value = 20.7
$ws.Range("V64").Value = 410.55
$ws.Range("W64").Value = 358.79
$ws.Range("Y64").Value = 228.82
$ws.Range("Z64").Value = 329.99
$ws.Range("V66").Value = 17.54
$ws.Range("W66").Value = 14.35
$ws.Range("Y66").Value = 7.78
$ws.Range("Z66").Value = 12.54
$ws.Range("V68").Value = 20.7
$ws.Range("V69").Value = 407.71
$ws.Range("W69").Value = 358.05
$ws.Range("Y69").Value = 227.64
$ws.Range("Z69").Value = 328.96
$ws.Range("V71").Value = 14.71
$ws.Range("W71").Value = 13.61
$ws.Range("Y71").Value = 6.6
$ws.Range("Z71").Value = 11.51
$ws.Range("V73").Value = 20.7
$ws.Range("V74").Value = 402.16
$ws.Range("W74").Value = 354
$ws.Range("Y74").Value = 225.78
$ws.Range("Z74").Value = 326.59
$ws.Range("V76").Value = 9.15
$ws.Range("W76").Value = 9.56
$ws.Range("Y76").Value = 4.74
$ws.Range("Z76").Value = 9.140000000000001
$ws.Range("V78").Value = 20.7
$ws.Range("V79").Value = 393.01
$ws.Range("W79").Value = 344.44
$ws.Range("Y79").Value = 221.04
$ws.Range("Z79").Value = 317.45
$ws.Range("V83").Value = 20.7
$ws.Range("V84").Value = 378.69
$ws.Range("W84").Value = 330.24
$ws.Range("X84").Value = 219.06
$ws.Range("Y84").Value = 224.63
$ws.Range("Z84").Value = 309.41
$ws.Range("V86").Value = -14.32
$ws.Range("W86").Value = -14.2
$ws.Range("X86").Value = -5.26
$ws.Range("Y86").Value = 3.59
$ws.Range("Z86").Value = -8.039999999999999
$ws.Range("V88").Value = 20.7
$ws.Range("V89").Value = 357.63
$ws.Range("W89").Value = 313.41
$ws.Range("Y89").Value = 204.29
$ws.Range("Z89").Value = 294.48
$ws.Range("V91").Value = -35.38
$ws.Range("W91").Value = -31.03
$ws.Range("Y91").Value = -16.75
$ws.Range("Z91").Value = -22.97
$ws.Range("V93").Value = 20.7
